$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 336.83334
$ws.Range("I28").Value = 304.2
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 304.2
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = 180.8
$ws.Range("N28").Value = -1470
# Row 112
$ws.Range("H112").Value = 2399.6667
$ws.Range("J112").Value = 3142.8572
$ws.Range("L112").Value = 9428.571599999999
$ws.Range("N112").Value = -11644.5716
# Row 113
$ws.Range("H113").Value = 8599.5
$ws.Range("I113").Value = 7199
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 7199
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -3945
$ws.Range("N113").Value = -16508
# Row 137
$ws.Range("H137").Value = 1868.72
$ws.Range("I137").Value = 1364.091
$ws.Range("K137").Value = 4092.273
$ws.Range("M137").Value = -1542.273

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1624.8077
$ws.Range("I32").Value = 1415.1915
$ws.Range("J32").Value = 3595.2
$ws.Range("K32").Value = 1415.1915
$ws.Range("L32").Value = 3595.2
$ws.Range("M32").Value = -1128.1915
$ws.Range("N32").Value = -4169.2
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
# Row 61
$ws.Range("H61").Value = 2590.9
$ws.Range("I61").Value = 2875.625
$ws.Range("J61").Value = 1452
$ws.Range("K61").Value = 2875.625
$ws.Range("L61").Value = 1452
$ws.Range("M61").Value = -2663.625
$ws.Range("N61").Value = -1876
# Row 92
$ws.Range("H92").Value = 28000
$ws.Range("J92").Value = 28000
$ws.Range("L92").Value = 28000
$ws.Range("N92").Value = -32992
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
# Row 110
$ws.Range("H110").Value = 447.5
$ws.Range("J110").Value = 447.5
$ws.Range("L110").Value = 447.5
$ws.Range("N110").Value = -4537.5
# Row 122
$ws.Range("H122").Value = 2449.1667
$ws.Range("I122").Value = 1939.2
$ws.Range("K122").Value = 5817.6
$ws.Range("M122").Value = -3367.6
# Row 132
$ws.Range("H132").Value = 2186.3333
$ws.Range("I132").Value = 2186.3333
$ws.Range("K132").Value = 6558.999899999999
$ws.Range("M132").Value = -4028.999899999999
# Row 136
$ws.Range("H136").Value = 2590.9
$ws.Range("I136").Value = 2875.625
$ws.Range("J136").Value = 1452
$ws.Range("K136").Value = 8626.875
$ws.Range("L136").Value = 4356
$ws.Range("M136").Value = -6076.875
$ws.Range("N136").Value = -9456

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2542.2
$ws.Range("I94").Value = 2542.2
$ws.Range("K94").Value = 2542.2
$ws.Range("M94").Value = -2091.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 48.285713
$ws.Range("I7").Value = 48.285713
$ws.Range("K7").Value = 48.285713
$ws.Range("M7").Value = 64.714287
# Row 60
$ws.Range("H60").Value = 24885.666
$ws.Range("J60").Value = 24995.666
$ws.Range("L60").Value = 24995.666
$ws.Range("N60").Value = -26017.666
# Row 99
$ws.Range("H99").Value = 4614.6665
$ws.Range("I99").Value = 3377.6667
$ws.Range("K99").Value = 3377.6667
$ws.Range("M99").Value = -1879.6667
# Row 105
$ws.Range("H105").Value = 1750
$ws.Range("I105").Value = 1750
$ws.Range("K105").Value = 1750
$ws.Range("M105").Value = -3
# Row 122
$ws.Range("H122").Value = 1564.5
$ws.Range("J122").Value = 1244.5
$ws.Range("L122").Value = 3733.5
$ws.Range("N122").Value = -8633.5
# Row 125
$ws.Range("H125").Value = 25000
$ws.Range("J125").Value = 25000
$ws.Range("L125").Value = 25000
$ws.Range("N125").Value = -29920
# Row 126
$ws.Range("H126").Value = 4614.6665
$ws.Range("I126").Value = 3377.6667
$ws.Range("K126").Value = 10133.0001
$ws.Range("M126").Value = -7663.000100000001
# Row 129
$ws.Range("H129").Value = 58000
$ws.Range("J129").Value = 58000
$ws.Range("L129").Value = 58000
$ws.Range("N129").Value = -68000
# Row 131
$ws.Range("H131").Value = 49999.5
$ws.Range("J131").Value = 49999.5
$ws.Range("L131").Value = 49999.5
$ws.Range("N131").Value = -60079.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 700
$ws.Range("J131").Value = 700
$ws.Range("L131").Value = 2100
$ws.Range("N131").Value = -12180

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 97.588234
$ws.Range("I2").Value = 91.454544
$ws.Range("J2").Value = 108.833336
$ws.Range("K2").Value = 91.454544
$ws.Range("L2").Value = 108.833336
$ws.Range("M2").Value = 21.545456
$ws.Range("N2").Value = -334.833336
# Row 62
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372
# Row 65
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864
# Row 102
$ws.Range("H102").Value = 6265
$ws.Range("J102").Value = 16148.5
$ws.Range("L102").Value = 16148.5
$ws.Range("N102").Value = -19392.5
# Row 132
$ws.Range("H132").Value = 2272.3572
$ws.Range("I132").Value = 2365.182
$ws.Range("K132").Value = 7095.545999999999
$ws.Range("M132").Value = -4565.545999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 894.375
$ws.Range("I16").Value = 808.1429000000001
$ws.Range("K16").Value = 808.1429000000001
$ws.Range("M16").Value = -638.1429000000001
# Row 22
$ws.Range("H22").Value = 985
$ws.Range("I22").Value = 648.5
$ws.Range("K22").Value = 648.5
$ws.Range("M22").Value = -353.5
# Row 27
$ws.Range("H27").Value = 985
$ws.Range("I27").Value = 648.5
$ws.Range("K27").Value = 648.5
$ws.Range("M27").Value = -541.5
# Row 46
$ws.Range("H46").Value = 2555.4443
$ws.Range("J46").Value = 2500
$ws.Range("L46").Value = 2500
$ws.Range("N46").Value = -2876
# Row 55
$ws.Range("H55").Value = 216.29411
$ws.Range("I55").Value = 241.75
$ws.Range("K55").Value = 241.75
$ws.Range("M55").Value = -68.75
# Row 61
$ws.Range("H61").Value = 4001.8333
$ws.Range("I61").Value = 4775.3335
$ws.Range("K61").Value = 4775.3335
$ws.Range("M61").Value = -4573.3335
# Row 63
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31498
# Row 66
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -97488
# Row 113
$ws.Range("H113").Value = 4001.8333
$ws.Range("I113").Value = 4775.3335
$ws.Range("K113").Value = 4775.3335
$ws.Range("M113").Value = -2605.3335
# Row 118
$ws.Range("H118").Value = 40895
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
# Row 122
$ws.Range("H122").Value = 3646.1428
$ws.Range("I122").Value = 3503.8333
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 10511.4999
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -8061.499899999999
$ws.Range("N122").Value = -18400
# Row 132
$ws.Range("H132").Value = 13738.444
$ws.Range("I132").Value = 9632.333000000001
$ws.Range("K132").Value = 28896.999
$ws.Range("M132").Value = -26366.999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 1483.4445
$ws.Range("I4").Value = 525
$ws.Range("J4").Value = 1757.2858
$ws.Range("K4").Value = 525
$ws.Range("L4").Value = 1757.2858
$ws.Range("M4").Value = -412
$ws.Range("N4").Value = -1983.2858
# Row 101
$ws.Range("H101").Value = 23750
$ws.Range("J101").Value = 23750
$ws.Range("L101").Value = 23750
$ws.Range("N101").Value = -30240
# Row 113
$ws.Range("H113").Value = 224.42857
$ws.Range("I113").Value = 99
$ws.Range("K113").Value = 297
$ws.Range("M113").Value = 1873
# Row 116
$ws.Range("H116").Value = 28624.75
$ws.Range("J116").Value = 28624.75
$ws.Range("L116").Value = 28624.75
$ws.Range("N116").Value = -37802.75
# Row 122
$ws.Range("H122").Value = 3254.9524
$ws.Range("I122").Value = 3078
$ws.Range("K122").Value = 9234
$ws.Range("M122").Value = -6784
# Row 132
$ws.Range("H132").Value = 2521.1
$ws.Range("I132").Value = 2521.1
$ws.Range("K132").Value = 7563.299999999999
$ws.Range("M132").Value = -5033.299999999999
# Row 136
$ws.Range("H136").Value = 4236.2666
$ws.Range("I136").Value = 4418.1113
$ws.Range("K136").Value = 13254.3339
$ws.Range("M136").Value = -10704.3339
